# Refresh currentAveragePrice / LevePrice / LeveProfit figures for the Leve tables
# (H-N columns) across all eight job sheets, per the latest market-board pull
# from the scheduled pricing runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 2679.9143
$ws.Range("I33").Value = 2817.4482
$ws.Range("K33").Value = 2817.4482
$ws.Range("M33").Value = -2588.4482

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 58824130
$ws.Range("I98").Value = 58824130
$ws.Range("K98").Value = 58824130
$ws.Range("M98").Value = -58822632

# Row 105: Ultimate Official Strategy Guide / Gazelleskin Codex
$ws.Range("H105").Value = 57495
$ws.Range("J105").Value = 57495
$ws.Range("L105").Value = 57495
$ws.Range("N105").Value = -64483

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1122.9215
$ws.Range("J112").Value = 1210.659
$ws.Range("L112").Value = 3631.977
$ws.Range("N112").Value = -5847.977000000001

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 58824130
$ws.Range("I122").Value = 58824130
$ws.Range("K122").Value = 176472390
$ws.Range("M122").Value = -176469940

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3056.2593
$ws.Range("I137").Value = 2090.0264
$ws.Range("J137").Value = 5351.0625
$ws.Range("K137").Value = 6270.0792
$ws.Range("L137").Value = 16053.1875
$ws.Range("M137").Value = -3720.0792
$ws.Range("N137").Value = -21153.1875

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 31: I Was a Teenage Wailer / Iron Alembic
$ws.Range("H31").Value = 14984.333
$ws.Range("I31").Value = 2857.375
$ws.Range("K31").Value = 2857.375
$ws.Range("M31").Value = -2563.375

# Row 46: Get Me the Usual / Heavy Steel Flanchard
$ws.Range("H46").Value = 6263.75
$ws.Range("J46").Value = 6263.75
$ws.Range("L46").Value = 6263.75
$ws.Range("N46").Value = -6901.75

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 6346114.5
$ws.Range("I74").Value = 8623104
$ws.Range("K74").Value = 8623104
$ws.Range("M74").Value = -8622230

# Row 76: Sometimes the South Wins / Titanium Mail of Fending
$ws.Range("H76").Value = 24499
$ws.Range("J76").Value = 24499
$ws.Range("L76").Value = 24499
$ws.Range("N76").Value = -25175

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 6346114.5
$ws.Range("I77").Value = 8623104
$ws.Range("K77").Value = 43115520
$ws.Range("M77").Value = -43111152

# Row 79: The Thriller of Autumn (L) / Titanium Mail of Fending
$ws.Range("H79").Value = 24499
$ws.Range("J79").Value = 24499
$ws.Range("L79").Value = 24499
$ws.Range("N79").Value = -26839

# Row 94: Setting the Stage / High Steel Helm of Maiming
$ws.Range("H94").Value = 47323
$ws.Range("J94").Value = 47323
$ws.Range("L94").Value = 47323
$ws.Range("N94").Value = -49125

# Row 106: Heads Will Roll / Molybdenum Headgear of Maiming
$ws.Range("H106").Value = 50783.332
$ws.Range("J106").Value = 50783.332
$ws.Range("L106").Value = 50783.332
$ws.Range("N106").Value = -53307.332

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 5951.9375
$ws.Range("I132").Value = 2257.9
$ws.Range("J132").Value = 12108.667
$ws.Range("K132").Value = 6773.700000000001
$ws.Range("L132").Value = 36326.001
$ws.Range("M132").Value = -4243.700000000001
$ws.Range("N132").Value = -41386.001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 132: Always Be Prepaired / Mountain Chromite Twinfangs
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 42502470
$ws.Range("I134").Value = 2139
$ws.Range("K134").Value = 6417
$ws.Range("M134").Value = -3882

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 19: Shielding Sales / Square Ash Shield
$ws.Range("H19").Value = 2854.077
$ws.Range("I19").Value = 1680.4667
$ws.Range("J19").Value = 4454.4546
$ws.Range("K19").Value = 1680.4667
$ws.Range("L19").Value = 4454.4546
$ws.Range("M19").Value = -1510.4667
$ws.Range("N19").Value = -4794.4546

# Row 24: What You Need / Square Ash Shield
$ws.Range("H24").Value = 2854.077
$ws.Range("I24").Value = 1680.4667
$ws.Range("J24").Value = 4454.4546
$ws.Range("K24").Value = 1680.4667
$ws.Range("L24").Value = 4454.4546
$ws.Range("M24").Value = -1510.4667
$ws.Range("N24").Value = -4794.4546

# Row 28: Militia on My Mind / Iron Lance
$ws.Range("H28").Value = 8824.5
$ws.Range("J28").Value = 8824.5
$ws.Range("L28").Value = 8824.5
$ws.Range("N28").Value = -9314.5

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 657664.8
$ws.Range("I31").Value = 11761.211
$ws.Range("K31").Value = 11761.211
$ws.Range("M31").Value = -11466.211

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 657664.8
$ws.Range("I34").Value = 11761.211
$ws.Range("K34").Value = 11761.211
$ws.Range("M34").Value = -11559.211

# Row 92: Walk the Walk / Beech Rod
$ws.Range("H92").Value = 95246.5
$ws.Range("J92").Value = 95246.5
$ws.Range("L92").Value = 95246.5
$ws.Range("N92").Value = -100238.5

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2720.186
$ws.Range("I132").Value = 2793.8108
$ws.Range("J132").Value = 2266.1667
$ws.Range("K132").Value = 8381.432400000002
$ws.Range("L132").Value = 6798.500100000001
$ws.Range("M132").Value = -5851.432400000002
$ws.Range("N132").Value = -11858.5001

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 264550.34
$ws.Range("I134").Value = 287058.8
$ws.Range("J134").Value = 1951.6666
$ws.Range("K134").Value = 861176.3999999999
$ws.Range("L134").Value = 5854.9998
$ws.Range("M134").Value = -858641.3999999999
$ws.Range("N134").Value = -10924.9998

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 69: Loving That Muffin Top / Ishgardian Muffin
$ws.Range("H69").Value = 3074.3333
$ws.Range("J69").Value = 3189.2
$ws.Range("L69").Value = 9567.599999999999
$ws.Range("N69").Value = -11189.6

# Row 72: Muffin of the Morn (L) / Ishgardian Muffin
$ws.Range("H72").Value = 3074.3333
$ws.Range("J72").Value = 3189.2
$ws.Range("L72").Value = 28702.8
$ws.Range("N72").Value = -36814.8

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 3413.4707
$ws.Range("I80").Value = 2840.7144
$ws.Range("K80").Value = 8522.143199999999
$ws.Range("M80").Value = -7586.143199999999

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 3413.4707
$ws.Range("I83").Value = 2840.7144
$ws.Range("K83").Value = 25566.4296
$ws.Range("M83").Value = -20886.4296

# Row 87: Soup That Eats Like a Knight / Clam Chowder
$ws.Range("H87").Value = 2075
$ws.Range("I87").Value = 2075
$ws.Range("K87").Value = 6225
$ws.Range("M87").Value = -4977

# Row 90: Like Ma Used to Make (L) / Clam Chowder
$ws.Range("H90").Value = 2075
$ws.Range("I90").Value = 2075
$ws.Range("K90").Value = 18675
$ws.Range("M90").Value = -12435

# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 1769.091
$ws.Range("I121").Value = 482.25
$ws.Range("J121").Value = 2504.4285
$ws.Range("K121").Value = 1446.75
$ws.Range("L121").Value = 7513.2855
$ws.Range("M121").Value = -136.75
$ws.Range("N121").Value = -10133.2855

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 4825.077
$ws.Range("I131").Value = 8833.333000000001
$ws.Range("J131").Value = 4302.2607
$ws.Range("K131").Value = 26499.999
$ws.Range("L131").Value = 12906.7821
$ws.Range("M131").Value = -21459.999
$ws.Range("N131").Value = -22986.7821

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 93: One Ring Circus / Triphane Ring of Slaying
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744

# Row 95: Chain of Command / Koppranickel Temple Chain
$ws.Range("H95").Value = 53463
$ws.Range("J95").Value = 53463
$ws.Range("L95").Value = 53463
$ws.Range("N95").Value = -58955

# Row 109: You're My Wonderhall / Hematite Earrings of Healing
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 2137
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 71439390
$ws.Range("I132").Value = 166669120
$ws.Range("K132").Value = 500007360
$ws.Range("M132").Value = -500004830

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 39: Quality over Quantity / Boarskin Himantes
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Row 81: I Need Your Glove Tonight / Dragonskin Gloves of Healing
$ws.Range("H81").Value = 42499
$ws.Range("J81").Value = 42499
$ws.Range("L81").Value = 42499
$ws.Range("N81").Value = -44495

# Row 84: Halonic Drake Handlers (L) / Dragonskin Gloves of Healing
$ws.Range("H84").Value = 42499
$ws.Range("J84").Value = 42499
$ws.Range("L84").Value = 127497
$ws.Range("N84").Value = -137481

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 66234.45
$ws.Range("I136").Value = 10971.272
$ws.Range("K136").Value = 32913.81600000001
$ws.Range("M136").Value = -30363.81600000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 104: Brimming with Confidence / Twinsilk Turban of Aiming
$ws.Range("H104").Value = 17871.666
$ws.Range("I104").Value = 24665
$ws.Range("K104").Value = 24665
$ws.Range("M104").Value = -21171

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1581.1428
$ws.Range("I126").Value = 1581.1428
$ws.Range("K126").Value = 4743.428400000001
$ws.Range("M126").Value = -2273.428400000001
